$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 33
$ws1.Range("F5").Value = 27
$ws1.Range("F6").Value = 558
$ws1.Range("F7").Value = 1742
$ws1.Range("F11").Value = 1888
$ws1.Range("F13").Value = 122
$ws1.Range("F15").Value = 7
$ws1.Range("F23").Value = 1018
$ws1.Range("F25").Value = 317

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 33
$ws4.Range("F5").Value = 27
$ws4.Range("F6").Value = 558
$ws4.Range("F7").Value = 1742
$ws4.Range("F12").Value = 1888
$ws4.Range("F14").Value = 122
$ws4.Range("F16").Value = 7
$ws4.Range("F24").Value = 1018
$ws4.Range("F26").Value = 317
